$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume table refresh (row data + Bittensor/TheGraph reorder).
# Cells whose new text is a plain decimal number (e.g. "12.60", "0.649") get a
# leading apostrophe so COM stores them as text -- matching the source file,
# where every Price/Volume(1h) cell is an explicit string, not a number.

$ws.Range("D2").Value = '69.592.02'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '3.493.01'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = "'612.09"
$ws.Range("E5").Value = '  +5.29%  '
$ws.Range("D6").Value = "'189.42"
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("D10").Value = "'0.649"
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("D11").Value = "'53.05"
$ws.Range("E11").Value = '  -3.08%  '
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("E12").Value = '  -4.07%  '
$ws.Range("D13").Value = "'9.49"
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '4.056.22'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("D15").Value = "'599.06"
$ws.Range("E15").Value = '  +4.64%  '
$ws.Range("D16").Value = '69.621.26'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = "'19.01"
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = "'12.60"
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '3.491.80'
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = "'17.21"
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("D23").Value = "'105.22"
$ws.Range("E23").Value = '  +12.10%  '
$ws.Range("E24").Value = '  +4.24%  '
$ws.Range("D25").Value = "'4.67"
$ws.Range("E25").Value = '  +2.77%  '
$ws.Range("D26").Value = "'3.04"
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("D28").Value = "'9.72"
$ws.Range("E28").Value = '  +4.80%  '
$ws.Range("D29").Value = "'33.29"
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Value = "'4.16"
$ws.Range("E31").Value = '  +11.08%  '
$ws.Range("D32").Value = "'12.51"
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("E35").Value = '  -5.57%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  +7.03%  '
$ws.Range("D38").Value = '3.626.15'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = "'0.394"
$ws.Range("E39").Value = '  -4.27%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = "'509.66"
$ws.Range("E40").Value = '  -6.15%  '
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("E42").Value = '  -3.57%  '
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = "'0.0462"
$ws.Range("E44").Value = '  -1.92%  '
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("E46").Value = '  +2.29%  '
$ws.Range("E47").Value = '  -4.53%  '
$ws.Range("D48").Value = "'8.77"
$ws.Range("E48").Value = '  -6.09%  '
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = "'132.01"
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").Value = "'1.35"
$ws.Range("E51").Value = '  -7.92%  '
